$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update vehicle identifier data (Patente/Motor/Chasis for rows 2-3) ---
$ws.Range("H2").Value = "ZZZ499"
$ws.Range("I2").Value = "ABC0987AX295"
$ws.Range("J2").Value = "MMAA09XFGS290"
$ws.Range("H3").Value = "ZZZ500"
$ws.Range("I3").Value = "ABC0987AX296"
$ws.Range("J3").Value = "MMAA09XFGS291"

# --- Split conditional formatting: H2:J3 becomes its own duplicate-values rule ---
$ws.Range("H2:J3").FormatConditions.Item(1).Delete()

$fc1 = $ws.Range("H6:J14").FormatConditions.AddUniqueValues()
$fc1.DupeUnique = 1
$fc1.Font.Color = 13408780
$fc1.Interior.Color = 13551615

$fc2 = $ws.Range("H2:J3").FormatConditions.AddUniqueValues()
$fc2.DupeUnique = 1
$fc2.Font.Color = 13408780
$fc2.Interior.Color = 13551615

# --- Move active selection to J8 ---
$ws.Range("J8").Select()
